{"js": "// Replace each two-digit-by-two-digit multiplication expression in the\n// worksheet table with its updated counterpart, per the commit's diff.\n// Every \"<old>=\" string in the document is unique, so a direct\n// search-and-replace per pair is unambiguous and safe.\nconst replacements = [\n  [\"40\u00d721=\", \"89\u00d767=\"],\n  [\"12\u00d785=\", \"16\u00d768=\"],\n  [\"11\u00d753=\", \"80\u00d727=\"],\n  [\"44\u00d788=\", \"24\u00d725=\"],\n  [\"74\u00d711=\", \"87\u00d762=\"],\n  [\"50\u00d745=\", \"29\u00d719=\"],\n  [\"71\u00d757=\", \"47\u00d767=\"],\n  [\"34\u00d778=\", \"20\u00d785=\"],\n  [\"33\u00d715=\", \"96\u00d763=\"],\n  [\"87\u00d776=\", \"11\u00d724=\"],\n  [\"65\u00d737=\", \"62\u00d763=\"],\n  [\"32\u00d792=\", \"66\u00d793=\"],\n  [\"55\u00d797=\", \"45\u00d731=\"],\n  [\"25\u00d731=\", \"77\u00d739=\"],\n  [\"15\u00d759=\", \"39\u00d736=\"],\n  [\"14\u00d796=\", \"84\u00d777=\"],\n  [\"64\u00d741=\", \"50\u00d793=\"],\n  [\"60\u00d737=\", \"83\u00d711=\"],\n  [\"28\u00d786=\", \"27\u00d734=\"],\n  [\"96\u00d722=\", \"20\u00d740=\"],\n  [\"95\u00d732=\", \"23\u00d716=\"],\n  [\"68\u00d797=\", \"15\u00d750=\"],\n  [\"13\u00d737=\", \"34\u00d734=\"],\n  [\"57\u00d754=\", \"78\u00d798=\"],\n  [\"86\u00d793=\", \"34\u00d798=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication expression in the\n# worksheet table with its updated counterpart, per the commit's diff.\n# Every \"<old>=\" string in the document is unique, so a direct\n# Find/Replace per pair is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"40\u00d721=\", \"89\u00d767=\"),\n    @(\"12\u00d785=\", \"16\u00d768=\"),\n    @(\"11\u00d753=\", \"80\u00d727=\"),\n    @(\"44\u00d788=\", \"24\u00d725=\"),\n    @(\"74\u00d711=\", \"87\u00d762=\"),\n    @(\"50\u00d745=\", \"29\u00d719=\"),\n    @(\"71\u00d757=\", \"47\u00d767=\"),\n    @(\"34\u00d778=\", \"20\u00d785=\"),\n    @(\"33\u00d715=\", \"96\u00d763=\"),\n    @(\"87\u00d776=\", \"11\u00d724=\"),\n    @(\"65\u00d737=\", \"62\u00d763=\"),\n    @(\"32\u00d792=\", \"66\u00d793=\"),\n    @(\"55\u00d797=\", \"45\u00d731=\"),\n    @(\"25\u00d731=\", \"77\u00d739=\"),\n    @(\"15\u00d759=\", \"39\u00d736=\"),\n    @(\"14\u00d796=\", \"84\u00d777=\"),\n    @(\"64\u00d741=\", \"50\u00d793=\"),\n    @(\"60\u00d737=\", \"83\u00d711=\"),\n    @(\"28\u00d786=\", \"27\u00d734=\"),\n    @(\"96\u00d722=\", \"20\u00d740=\"),\n    @(\"95\u00d732=\", \"23\u00d716=\"),\n    @(\"68\u00d797=\", \"15\u00d750=\"),\n    @(\"13\u00d737=\", \"34\u00d734=\"),\n    @(\"57\u00d754=\", \"78\u00d798=\"),\n    @(\"86\u00d793=\", \"34\u00d798=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
